# Interdiff between v2 and v3 - applies the meaningful content edits made to
# the Logic Component Sequence Diagram slide (slide 1):
#   - The "Parser" lifeline box is renamed/re-wrapped to ":Address" / "BookParser"
#     (two lines) and resized/repositioned slightly.
#   - The "p" + "arse(\u201c1\u201d)" text box runs are consolidated into a single
#     "parse(\u201c1\u201d)" run.

function Get-ShapeById {
    param($Slide, [int]$Id)
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $shp = $Slide.Shapes.Item($i)
        if ($shp.Id -eq $Id) {
            return $shp
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 16 ("Rectangle 62"): ":Parser" -> ":Address" / "BookParser" ---
$parserShape = Get-ShapeById $s 16

# Reposition / resize the box slightly to fit the new two-line label.
$parserShape.Left = 216
$parserShape.Top = 171.59913385826772
$parserShape.Width = 96
$parserShape.Height = 36.825511811023624

# Force a full run rebuild (rather than an in-place text patch) so the two
# paragraphs come out as clean single runs, matching how PowerPoint rewrites
# the paragraph/run structure when the text is retyped.
$parserShape.TextFrame.TextRange.Text = "X"
$parserShape.TextFrame.TextRange.Text = ":Address" + [char]13 + "BookParser"

# --- Shape 79 ("TextBox 78"): merge "p" + "arse(...)" runs into "parse(...)" ---
$parseShape = Get-ShapeById $s 79
$quote1 = [char]0x201C
$quote2 = [char]0x201D
$parseShape.TextFrame.TextRange.Text = "X"
$parseShape.TextFrame.TextRange.Text = "parse(" + $quote1 + "1" + $quote2 + ")"
